$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($r in 2..3) {
    $ws.Cells.Item($r, 4).Value = 0.127   # D - historical_growth_revenue_last_5_years
    $ws.Range("E$r").ClearContents()       # E - historical_growth_net_income_last_5_years (removed)

    $ws.Range("G$r").Value = 0.1563234089046708
    $ws.Range("H$r").Value = 0.1012564872985523
    $ws.Range("I$r").Value = 0.1646064496752328
    $ws.Range("J$r").Value = 0.1646064496752328
    $ws.Range("K$r").Value = 670.4
    $ws.Range("L$r").Value = 0.1831193662933625
    $ws.Range("M$r").Value = 233.6
    $ws.Range("N$r").Value = 0.01459689815914118
    $ws.Range("O$r").Value = 0.3484486873508353
    $ws.Range("P$r").Value = 146.7
    $ws.Range("Q$r").Value = 0.009166802054563404
    $ws.Range("R$r").Value = 0.2188245823389021
    $ws.Range("S$r").Value = 86.90000000000001
    $ws.Range("T$r").Value = 0.3720034246575343
    $ws.Range("U$r").Value = 917.2
    $ws.Range("V$r").Value = 0.05731282102553208
    $ws.Range("W$r").Value = 0.5742184154175589
    $ws.Range("X$r").Value = 0.06768431463869683
    $ws.Range("Y$r").Value = 0.5065341007788621
    $ws.Range("Z$r").Value = 13.06664953851902
    $ws.Range("AA$r").Value = 2.150854789686135
    $ws.Range("AB$r").Value = 0.06758545600414918
    $ws.Range("AC$r").Value = 2.083269333681986
    $ws.Range("AD$r").Value = 0
    $ws.Range("AE$r").Value = 30.87893869486439
    $ws.Range("AF$r").Value = 30.87893869486439
    $ws.Range("AG$r").Value = -886.3210613051357
    $ws.Range("AH$r").Value = 0.001925807753059947
    $ws.Range("AI$r").Value = 0.01777840585600492
    $ws.Range("AJ$r").Value = -0.05863044473733868
    $ws.Range("AK$r").Value = -1.081302714348599
    $ws.Range("AL$r").Value = 0
    $ws.Range("AM$r").Value = -5.81
    $ws.Range("AN$r").Value = 0
    $ws.Range("AP$r").Value = -1.303030081307168
    $ws.Range("AQ$r").Value = -102.7710843373494
}
